# The invoice data table on the active sheet had its rows re-keyed against
# an updated source tree / DB schema: row 2 (HYD8-1345061 / CLICKTECH) and
# row 6 (FAAOZK2200001887 / Sane Retails) swap places, the verbose Amazon
# "Laptop Model" description is replaced with the short canonical model
# name "MacBook Air" everywhere it appears, and a couple of stray
# whitespace glitches in the incoming row (leading space in the invoice
# number, stray space in the serial number) are cleaned up.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: becomes the former row-6 (FAAOZK2200001887) data, cleaned up ---
# B2/C2 ("02-04-2021", "04-04-2021") look like valid MM-DD-YYYY dates to
# Excel's smart-entry parser, so force the cells to Text format first -
# otherwise they'd silently turn into date serials instead of the literal
# strings the source data uses (cf. untouched text dates like B3).
$ws.Range("A2").Value = "FAAOZK2200001887"
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "02-04-2021"
$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "04-04-2021"
$ws.Range("D2").Value = "OD121420552428604000"
$ws.Range("E2").Value = "Sane Retails Private Limited"
$ws.Range("F2").Value = "MacBook Air"
$ws.Range("G2").Value = "i5 5th Gen"
$ws.Range("H2").Value = "8GB"
$ws.Range("I2").Value = "128GB SSD"
$ws.Range("J2").Value = ""
$ws.Range("K2").Value = ""
$ws.Range("L2").Value = "Mac OS Sierra"
$ws.Range("M2").Value = ""
$ws.Range("N2").Value = "SFVFDL811J1WK"
$ws.Range("O2").Value = 12
$ws.Range("P2").Value = 67640

# --- Rows 3-5: only the verbose laptop-model description is shortened ---
$ws.Range("F3").Value = "MacBook Air"
$ws.Range("F4").Value = "MacBook Air"
$ws.Range("F5").Value = "MacBook Air"

# --- Row 6: becomes the former row-2 (HYD8-1345061) data ---
$ws.Range("A6").Value = "HYD8-1345061"
$ws.Range("B6").Value = "26-09-2024"
$ws.Range("C6").Value = "26-09-2024"
$ws.Range("D6").Value = "403-0151201-6865928"
$ws.Range("E6").Value = "CLICKTECH RETAIL PRIVATE LIMITED"
$ws.Range("F6").Value = "MacBook Air"
$ws.Range("G6").Value = "M1"
$ws.Range("H6").Value = "8GB"
$ws.Range("I6").Value = "256GB SSD"
$ws.Range("J6").Value = "Space Grey"
$ws.Range("K6").Value = "13.3-inch"
$ws.Range("L6").Value = ""
$ws.Range("N6").Value = ""
$ws.Range("O6").Value = ""
$ws.Range("P6").Value = 44906.78
